# Weekly update: insert a new data row at row 13 (most recent week),
# pushing the existing rows 13..81 down to 14..82.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value  = 10
$ws.Cells.Item(13, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(13, 3).Value  = 'La Araucanía'
$ws.Cells.Item(13, 4).Value  = 44558
$ws.Cells.Item(13, 5).Value  = 9
$ws.Cells.Item(13, 6).Value  = 'Fruta'
$ws.Cells.Item(13, 7).Value  = 100101
$ws.Cells.Item(13, 8).Value  = 'Berries'
$ws.Cells.Item(13, 9).Value  = 100101001
$ws.Cells.Item(13, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(13, 11).Value = 'Sin especificar'
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 500
$ws.Cells.Item(13, 14).Value = 2200
$ws.Cells.Item(13, 15).Value = 2200
$ws.Cells.Item(13, 16).Value = 2200
$ws.Cells.Item(13, 17).Value = '$/kilo'
$ws.Cells.Item(13, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(13, 19).Value = 2200
$ws.Cells.Item(13, 20).Value = 1
